# Append the new bitcoin purchase row (2025-10-26 run) to the sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newRow = 54

# Column A holds the purchase date as plain text (e.g. "10/19/2025" on the
# row above), not a real Excel date. Force text formatting before writing
# the value so Excel doesn't auto-coerce the "MM/DD/YYYY" string into a
# date serial, then drop the explicit formatting again so the cell ends up
# with the same default (unstyled) look as the rest of the data rows.
$dateCell = $ws.Range("A" + $newRow)
$dateCell.NumberFormat = "@"
$dateCell.Value = "10/26/2025"
$dateCell.ClearFormats()

$ws.Range("B" + $newRow).Value = 0.000445999999999995
$ws.Range("C" + $newRow).Value = 112107.6233183869
$ws.Range("D" + $newRow).Value = 50
